# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# 8f104349-c4d2-4df8-be52-d8076a42e2d6.md file is now "Ready for handoff",
# including an error detail noting the handback file version mismatch.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80eab20ee784c6e10fc7d2095f7890c70020abf5/e2e/8f104349-c4d2-4df8-be52-d8076a42e2d6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95b4f27b504030077da3efd4bff6afd8fcc30464/e2e/8f104349-c4d2-4df8-be52-d8076a42e2d6.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 12:52:29"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-04 12:52:25"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-04 12:52:29"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
